# Updates cryptos list values (Price column D, Volume(1h) column E)
# per the "Updated cryptos list" GitHub Actions commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value2 = "28.038.80"
$ws.Range("E2").Value2 = "  -2.01%  "
$ws.Range("D3").Value2 = "1.830.93"
$ws.Range("E3").Value2 = "  -1.08%  "
$ws.Range("E4").Value2 = "  -0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value2 = "323.71"
$ws.Range("E5").Value2 = "  -3.13%  "
$ws.Range("E6").Value2 = "  +0.00%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value2 = "0.4660"
$ws.Range("E7").Value2 = "  +0.34%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value2 = "0.3863"
$ws.Range("E8").Value2 = "  -1.40%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value2 = "0.07861"
$ws.Range("E9").Value2 = "  -0.70%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value2 = "0.9578"
$ws.Range("E10").Value2 = "  -2.86%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value2 = "21.88"
$ws.Range("E11").Value2 = "  -1.87%  "
$ws.Range("D12").Value2 = "1.852.01"
$ws.Range("E12").Value2 = "  -5.53%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value2 = "5.679"
$ws.Range("E13").Value2 = "  -2.92%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value2 = "6.900"
$ws.Range("E14").Value2 = "  -1.65%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value2 = "0.06851"
$ws.Range("E15").Value2 = "  -0.20%  "
$ws.Range("E16").Value2 = "  -0.84%  "
$ws.Range("E17").Value2 = "  +0.00%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value2 = "0.000009911"
$ws.Range("E18").Value2 = "  -1.59%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value2 = "16.58"
$ws.Range("E19").Value2 = "  -3.22%  "
$ws.Range("E20").Value2 = "  -0.02%  "
$ws.Range("D21").Value2 = "28.041.07"
$ws.Range("E21").Value2 = "  -2.11%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value2 = "5.317"
$ws.Range("E22").Value2 = "  -1.49%  "
$ws.Range("E23").Value2 = "  -3.06%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value2 = "2.094"
$ws.Range("E24").Value2 = "  -1.62%  "
$ws.Range("D25").Value2 = "2.113.06"
$ws.Range("E25").Value2 = "  -4.90%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value2 = "153.78"
$ws.Range("E26").Value2 = "  +0.47%  "
$ws.Range("E27").Value2 = "  -1.73%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value2 = "5.677"
$ws.Range("E28").Value2 = "  -7.07%  "
$ws.Range("E29").Value2 = "  -3.11%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value2 = "117.62"
$ws.Range("E30").Value2 = "  +0.07%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value2 = "0.09252"
$ws.Range("E31").Value2 = "  -1.74%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value2 = "0.9323"
$ws.Range("E32").Value2 = "  -4.85%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value2 = "5.265"
$ws.Range("E33").Value2 = "  -1.94%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value2 = "1.319"
$ws.Range("E34").Value2 = "  -2.40%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value2 = "3.295"
$ws.Range("E35").Value2 = "  -5.26%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value2 = "0.05831"
$ws.Range("E36").Value2 = "  -5.21%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value2 = "0.02117"
$ws.Range("E37").Value2 = "  -3.83%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value2 = "1.132"
$ws.Range("E38").Value2 = "  -2.68%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value2 = "7.797"
$ws.Range("E39").Value2 = "  +2.16%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value2 = "0.5586"
$ws.Range("E40").Value2 = "  -2.22%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value2 = "9.861"
$ws.Range("E41").Value2 = "  -2.80%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value2 = "0.1759"
$ws.Range("E42").Value2 = "  -2.20%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value2 = "0.07269"
$ws.Range("E43").Value2 = "  +1.58%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value2 = "11.66"
$ws.Range("E44").Value2 = "  -1.18%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value2 = "0.5260"
$ws.Range("E45").Value2 = "  -2.55%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value2 = "1.828"
$ws.Range("E48").Value2 = "  -4.32%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value2 = "112.32"
$ws.Range("E49").Value2 = "  -1.61%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value2 = "1.001"
$ws.Range("E50").Value2 = "  +0.03%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value2 = "2.319"
$ws.Range("E51").Value2 = "  +0.17%  "

# Rows 46 and 47 swap positions: WEMIXToken now ranks above RenderToken
$ws.Range("B46").Value2 = "WEMIXToken"
$ws.Range("C46").Value2 = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value2 = "1.131"
$ws.Range("E46").Value2 = "  -9.35%  "

$ws.Range("B47").Value2 = "RenderToken"
$ws.Range("C47").Value2 = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value2 = "2.119"
$ws.Range("E47").Value2 = "  -11.42%  "
